$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values for the full A1:C26 label/translation table ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOQ4056"
$ws.Range("C2").Value = "LOQ4056"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Química Analítica para Engenharia"
$ws.Range("C3").Value = " Química Analítica para Engenharia"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Analytical Chemistry for Engineering"
$ws.Range("C4").Value = "Analytical Chemistry for Engineering"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2012"
$ws.Range("C8").Value = "01/01/2012"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EB-4,EQD-5,EQN-5"
$ws.Range("C9").Value = "EB-4,EQD-5,EQN-5"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Gerais: - Mostrar a Química Analítica por via úmida como uma ciência que se propõe a determinar a composição qualitativa e quantitativa da matéria por meio de reações químicas específicas e observação crítica dos resultados, requerendo para isso observadores competentes tecnicamente, criativos e sensatos.`nEspecíficos: - Ao concluir o curso os alunos devem: interpretar adequadamente as técnicas e princípios inseridos nos textos de Química Analítica; saber manusear com precisão e eficiência a instrumentação analítica, produtos tóxicos, inflamáveis e cáusticos; compreender os diversos tipos de cálculos estequiométricos; preparar, aferir, conservar e usar adequadamente soluções padrões, bem como, o descarte adequadamente em função da toxicidade dos reagentes/produtos."
$ws.Range("C10").Value = "Gerais: - Mostrar a Química Analítica por via úmida como uma ciência que se propõe a determinar a composição qualitativa e quantitativa da matéria por meio de reações químicas específicas e observação crítica dos resultados, requerendo para isso observadores competentes tecnicamente, criativos e sensatos.`nEspecíficos: - Ao concluir o curso os alunos devem: interpretar adequadamente as técnicas e princípios inseridos nos textos de Química Analítica; saber manusear com precisão e eficiência a instrumentação analítica, produtos tóxicos, inflamáveis e cáusticos; compreender os diversos tipos de cálculos estequiométricos; preparar, aferir, conservar e usar adequadamente soluções padrões, bem como, o descarte adequadamente em função da toxicidade dos reagentes/produtos."

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "General- Show the Analytical Chemistry wet method as a science that objective to determine the qualitative and quantitative composition of matter through specific chemical reactions and critical observation of the results, requiring only technically competent observers, creative and sensible.Specifics: - By completing the course students should: properly interpret the principles and techniques inserted into in the texts of Analytical Chemistry; know accurately and efficiently handle the analytical instrumentation, toxic, flammable and caustic; understand the various types of stoichiometric calculations, prepare, measure, solutions properly maintain and use patterns, as well as the disposal appropriately depending on the toxicity of the reactants / products."
$ws.Range("C11").Value = "General- Show the Analytical Chemistry wet method as a science that objective to determine the qualitative and quantitative composition of matter through specific chemical reactions and critical observation of the results, requiring only technically competent observers, creative and sensible.Specifics: - By completing the course students should: properly interpret the principles and techniques inserted into in the texts of Analytical Chemistry; know accurately and efficiently handle the analytical instrumentation, toxic, flammable and caustic; understand the various types of stoichiometric calculations, prepare, measure, solutions properly maintain and use patterns, as well as the disposal appropriately depending on the toxicity of the reactants / products."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "5840601 - Hélcio José Izário Filho"
$ws.Range("C13").Value = "5840601 - Hélcio José Izário Filho"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "- Bases teóricas da química analítica; Introdução à análise qualitativa; Leis e teorias fundamentais; Análise qualitativa sistemática de cátions; Análise qualitativa de ânions.`n- Fundamentos da análise titrimétrica; Titrimetria por Neutralização; Titrimetria por Precipitação; Titrimetria por oxidação-redução: Permanganatometria e Tiossulfatometria; Titrimetria por Complexação."
$ws.Range("C14").Value = "- Bases teóricas da química analítica; Introdução à análise qualitativa; Leis e teorias fundamentais; Análise qualitativa sistemática de cátions; Análise qualitativa de ânions.`n- Fundamentos da análise titrimétrica; Titrimetria por Neutralização; Titrimetria por Precipitação; Titrimetria por oxidação-redução: Permanganatometria e Tiossulfatometria; Titrimetria por Complexação."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "- Theoretical bases of analytical chemistry, introduction to qualitative analysis; Laws and theories; systematic qualitative analysis of cations, anions qualitative analysis.- Fundamentals of analysis titrimetry: titrimetry by Neutralization; Precipitation, redox titrimetry by KMnO4 and Na2S2O3; titrimetry by complexation."
$ws.Range("C15").Value = "- Theoretical bases of analytical chemistry, introduction to qualitative analysis; Laws and theories; systematic qualitative analysis of cations, anions qualitative analysis.- Fundamentals of analysis titrimetry: titrimetry by Neutralization; Precipitation, redox titrimetry by KMnO4 and Na2S2O3; titrimetry by complexation."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "- Bases Teóricas da Análise Qualitativa: Equilíbrio Químico; Efeito do íon Comum; Produto Iônico da água; Concentração do íon H+; Soluções tampão; Hidrólise dos sais; Produto de solubilidade. Operações analíticas: Precipitação; Filtração, Centrifugação. Equipamentos para ensaios por via úmida. Limpeza da aparelhagem/vidrarias. Classificação analítica dos cátions e dos ânions. Análise Qualitativa Sistemática: separação e identificação dos cátions do 1º, 2º e 3º Grupos. Identificação dos ânions segundo Vogel.`n- Fundamentos de Análise Quantitativa - Titrimetria por Neutralização: fundamentos específicos. Preparação e aferição das soluções padrão ácidas e alcalinas; determinações alcalimétricas e acidimétricas. - Titrimetria por Precipitação: discussão geral da Argentimetria. Princípio da acão dos indicadores. Preparação e emprego da solução padrão de nitrato de prata. Sulfocianetometria. Discussão geral. Preparação, aferição e emprego do processo. - Titrimetria por Oxidação-redução: - Permanganatometria. Características gerais do processo. Preparação / aferição e emprego do processo. - Tiossulfatometria: Características gerais do método. Emprego dos processos titulométricos: direto, inverso, indireto e de retorno. - Complexometria: Estudo teórico da formação de complexos. Grupos de coordenação. Características gerais do método. Preparação, aferição e emprego da solução padrão de EDTAH2Na2. Uso de indicadores metalocrômicos. Determinação de metais bivalentes e trivalentes."
$ws.Range("C16").Value = "- Bases Teóricas da Análise Qualitativa: Equilíbrio Químico; Efeito do íon Comum; Produto Iônico da água; Concentração do íon H+; Soluções tampão; Hidrólise dos sais; Produto de solubilidade. Operações analíticas: Precipitação; Filtração, Centrifugação. Equipamentos para ensaios por via úmida. Limpeza da aparelhagem/vidrarias. Classificação analítica dos cátions e dos ânions. Análise Qualitativa Sistemática: separação e identificação dos cátions do 1º, 2º e 3º Grupos. Identificação dos ânions segundo Vogel.`n- Fundamentos de Análise Quantitativa - Titrimetria por Neutralização: fundamentos específicos. Preparação e aferição das soluções padrão ácidas e alcalinas; determinações alcalimétricas e acidimétricas. - Titrimetria por Precipitação: discussão geral da Argentimetria. Princípio da acão dos indicadores. Preparação e emprego da solução padrão de nitrato de prata. Sulfocianetometria. Discussão geral. Preparação, aferição e emprego do processo. - Titrimetria por Oxidação-redução: - Permanganatometria. Características gerais do processo. Preparação / aferição e emprego do processo. - Tiossulfatometria: Características gerais do método. Emprego dos processos titulométricos: direto, inverso, indireto e de retorno. - Complexometria: Estudo teórico da formação de complexos. Grupos de coordenação. Características gerais do método. Preparação, aferição e emprego da solução padrão de EDTAH2Na2. Uso de indicadores metalocrômicos. Determinação de metais bivalentes e trivalentes."

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "- Theoretical Bases of Qualitative Analysis: Chemical Equilibrium; Effect of Common Ion, Ionic product of water, the H+ ion concentration, buffer solutions, hydrolysis of salts; Product solubility. Analytical operations: Precipitation, filtration, centrifugation. Equipment for testing wet method. Cleaning of equipment / glassware. Analytical classification of cations and anions. Qualitative Systematic Review: separation and identification of cations in the 1st, 2, 3rd and 4th groups. Identification of anions by Vogel.- Foundations of Quantitative Analysis - Titrimetry by Neutralization: specific discussion; Preparation and measurement of acid and alkaline buffer solutions; alkalimetry and acidimetry determinations. - Precipitation titrimetry: general discussion of argentometry. Principle of action of indicator. Preparation and use of standard solution of silver nitrate. thiocyanatemetry. General discussion. Preparation, examination and employment process. - Oxidation - reduction titrimetry: -permanganatemetry. General characteristics of the process. Preparation / evaluation and employment process. - thiocyanatemetry: General characteristics of the method. Employment titrimetric processes: direct, inverse, indirect and return. - Complexometry: Theoretical study of the formation of complexes. Coordination groups. General characteristics of the method. Preparation, assessment and use of standard solution EDTAH2Na2. Use of complexion indicators. Determination of bivalent and trivalent metals."
$ws.Range("C17").Value = "- Theoretical Bases of Qualitative Analysis: Chemical Equilibrium; Effect of Common Ion, Ionic product of water, the H+ ion concentration, buffer solutions, hydrolysis of salts; Product solubility. Analytical operations: Precipitation, filtration, centrifugation. Equipment for testing wet method. Cleaning of equipment / glassware. Analytical classification of cations and anions. Qualitative Systematic Review: separation and identification of cations in the 1st, 2, 3rd and 4th groups. Identification of anions by Vogel.- Foundations of Quantitative Analysis - Titrimetry by Neutralization: specific discussion; Preparation and measurement of acid and alkaline buffer solutions; alkalimetry and acidimetry determinations. - Precipitation titrimetry: general discussion of argentometry. Principle of action of indicator. Preparation and use of standard solution of silver nitrate. thiocyanatemetry. General discussion. Preparation, examination and employment process. - Oxidation - reduction titrimetry: -permanganatemetry. General characteristics of the process. Preparation / evaluation and employment process. - thiocyanatemetry: General characteristics of the method. Employment titrimetric processes: direct, inverse, indirect and return. - Complexometry: Theoretical study of the formation of complexes. Coordination groups. General characteristics of the method. Preparation, assessment and use of standard solution EDTAH2Na2. Use of complexion indicators. Determination of bivalent and trivalent metals."

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Serão aplicadas, por bimestre, duas avaliações, sendo uma avaliação teórica (peso 0,6) e uma avaliação prática (peso 0,4)."
$ws.Range("C19").Value = "Serão aplicadas, por bimestre, duas avaliações, sendo uma avaliação teórica (peso 0,6) e uma avaliação prática (peso 0,4)."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "A composição da média P1 e P2 será calculado pelo valor da avaliação teórica x 0,6 mais o valor da avaliação prática x 0,4. A média final será a média aritmética da P1 e P2."
$ws.Range("C20").Value = "A composição da média P1 e P2 será calculado pelo valor da avaliação teórica x 0,6 mais o valor da avaliação prática x 0,4. A média final será a média aritmética da P1 e P2."

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Na semana da recuperação será dado uma aula teórica e uma avaliação teórica no valor de 10. A Nota final será a média entre a média final (P1 e P2) e a nota da recuperação."
$ws.Range("C21").Value = "Na semana da recuperação será dado uma aula teórica e uma avaliação teórica no valor de 10. A Nota final será a média entre a média final (P1 e P2) e a nota da recuperação."

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Bibliografia Básica: `n1) VOGEL, Arthur Israel. Química analítica qualitativa. Sao Paulo: Mestrejou, 1981.`n2) VOGEL, Arthur I. Análise química quantitativa/ G. H. Jeffery; J. Bassett; J. Mendham; R. C. Denney. Rio de Janeiro: Guanabara Koogan, 1992.`n3) BACCAN, Nivaldo; ANDRADE, João Carlos de; GODINHO, Oswaldo E.S.; BARONE, José Salvador. Química analítica quantitativa elementar. São Paulo: Edgard Blücher - Instituto Mauá de Tecnologia, 2005-2007.`n4) BACCAN, Nivaldo et al.  Introdução à semimicroanálise qualitativa. Campinas:Editora da UNICAMP, 1988.`nBibliografia Complementar:`nSKOOG, Douglas A. et al. Fundamentos da química analítica. São Paulo: Editora Thomson Learning, 2006-9. `nALEXEYEV, V. Análise Qualitativa. Porto: Editora Lopes da Silva, 1982. `nHARRIS, Daniel C. Análise Química Quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2005."
$ws.Range("C22").Value = "Bibliografia Básica: `n1) VOGEL, Arthur Israel. Química analítica qualitativa. Sao Paulo: Mestrejou, 1981.`n2) VOGEL, Arthur I. Análise química quantitativa/ G. H. Jeffery; J. Bassett; J. Mendham; R. C. Denney. Rio de Janeiro: Guanabara Koogan, 1992.`n3) BACCAN, Nivaldo; ANDRADE, João Carlos de; GODINHO, Oswaldo E.S.; BARONE, José Salvador. Química analítica quantitativa elementar. São Paulo: Edgard Blücher - Instituto Mauá de Tecnologia, 2005-2007.`n4) BACCAN, Nivaldo et al.  Introdução à semimicroanálise qualitativa. Campinas:Editora da UNICAMP, 1988.`nBibliografia Complementar:`nSKOOG, Douglas A. et al. Fundamentos da química analítica. São Paulo: Editora Thomson Learning, 2006-9. `nALEXEYEV, V. Análise Qualitativa. Porto: Editora Lopes da Silva, 1982. `nHARRIS, Daniel C. Análise Química Quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2005."

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"

$ws.Range("B25").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"

$ws.Range("B26").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("C26").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"

# --- Clear cells that must no longer exist ---
$ws.Range("A13").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# --- Reset rows that should have default (non-custom) height ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(23).AutoFit()

# --- Apply explicit custom row heights ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30

